$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Summary")
$ws2 = $wb.Worksheets.Item("Cards_telegram")

# ---- Sheet1 (Summary): update rows 2-4, add row 5 ----
# Row 2: Сибирь – Барыс
$ws1.Cells.Item(2, 1).Value = 1369
$ws1.Cells.Item(2, 2).Value = 45992.64583333334
$ws1.Cells.Item(2, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(2, 3).Value = "Сибирь"
$ws1.Cells.Item(2, 4).Value = "Барыс"
$ws1.Cells.Item(2, 5).Value = "Сибирь – Барыс"
$ws1.Cells.Item(2, 6).Value = 897836
$ws1.Cells.Item(2, 7).Value = "https://text.khl.ru/text/897836.html"
$ws1.Cells.Item(2, 8).Value = 0.846154
$ws1.Cells.Item(2, 9).Value = 1.854302
$ws1.Cells.Item(2, 10).Value = 2.700456
$ws1.Cells.Item(2, 11).Value = 23.063363
$ws1.Cells.Item(2, 12).Value = 29.982381
$ws1.Cells.Item(2, 13).Value = 53.045744
$ws1.Cells.Item(2, 14).Value = 0.423308
$ws1.Cells.Item(2, 15).Value = 0.17033
$ws1.Cells.Item(2, 16).Value = 0.405949
$ws1.Cells.Item(2, 17).Value = 2.362346093152031
$ws1.Cells.Item(2, 18).Value = 5.870956378794105
$ws1.Cells.Item(2, 19).Value = 2.463363624494702
$ws1.Cells.Item(2, 20).Value = 42.3308
$ws1.Cells.Item(2, 21).Value = 17.033
$ws1.Cells.Item(2, 22).Value = 40.5949
$ws1.Cells.Item(2, 23).Value = 0.319303
$ws1.Cells.Item(2, 24).Value = 0.680284
$ws1.Cells.Item(2, 25).Value = 1.469974304849151
$ws1.Cells.Item(2, 26).Value = 0.485923
$ws1.Cells.Item(2, 27).Value = 0.513664
$ws1.Cells.Item(2, 28).Value = 1.946797906802891
$ws1.Cells.Item(2, 29).Value = 0.645717
$ws1.Cells.Item(2, 30).Value = 0.35387
$ws1.Cells.Item(2, 31).Value = 2.825896515669596
$ws1.Cells.Item(2, 32).Value = 0.786022
$ws1.Cells.Item(2, 33).Value = 0.213978
$ws1.Cells.Item(2, 34).Value = 1.272229021579549
$ws1.Cells.Item(2, 35).Value = 0.554926
$ws1.Cells.Item(2, 36).Value = 0.445074
$ws1.Cells.Item(2, 37).Value = 1.802042074078345
$ws1.Cells.Item(2, 38).Value = 0.777357
$ws1.Cells.Item(2, 39).Value = 0.222643
$ws1.Cells.Item(2, 40).Value = 1.286410233650691
$ws1.Cells.Item(2, 41).Value = 0.542462
$ws1.Cells.Item(2, 42).Value = 0.457538
$ws1.Cells.Item(2, 43).Value = 1.843447098598612
$ws1.Cells.Item(2, 44).Value = 0.746927
$ws1.Cells.Item(2, 45).Value = 1.338818920724515
$ws1.Cells.Item(2, 46).Value = 0.732445
$ws1.Cells.Item(2, 47).Value = 1.365290226569913

# Row 3: Северсталь – Трактор
$ws1.Cells.Item(3, 1).Value = 1369
$ws1.Cells.Item(3, 2).Value = 45992.79166666666
$ws1.Cells.Item(3, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(3, 3).Value = "Северсталь"
$ws1.Cells.Item(3, 4).Value = "Трактор"
$ws1.Cells.Item(3, 5).Value = "Северсталь – Трактор"
$ws1.Cells.Item(3, 6).Value = 897835
$ws1.Cells.Item(3, 7).Value = "https://text.khl.ru/text/897835.html"
$ws1.Cells.Item(3, 8).Value = 1.46875
$ws1.Cells.Item(3, 9).Value = 3.9375
$ws1.Cells.Item(3, 10).Value = 5.40625
$ws1.Cells.Item(3, 11).Value = 23.963798
$ws1.Cells.Item(3, 12).Value = 35.737266
$ws1.Cells.Item(3, 13).Value = 59.701065
$ws1.Cells.Item(3, 14).Value = 0.427646
$ws1.Cells.Item(3, 15).Value = 0.179644
$ws1.Cells.Item(3, 16).Value = 0.392535
$ws1.Cells.Item(3, 17).Value = 2.338382681002511
$ws1.Cells.Item(3, 18).Value = 5.566564984079625
$ws1.Cells.Item(3, 19).Value = 2.547543531150088
$ws1.Cells.Item(3, 20).Value = 42.7646
$ws1.Cells.Item(3, 21).Value = 17.9644
$ws1.Cells.Item(3, 22).Value = 39.2535
$ws1.Cells.Item(3, 23).Value = 0.406673
$ws1.Cells.Item(3, 24).Value = 0.593153
$ws1.Cells.Item(3, 25).Value = 1.685905660091072
$ws1.Cells.Item(3, 26).Value = 0.581479
$ws1.Cells.Item(3, 27).Value = 0.418346
$ws1.Cells.Item(3, 28).Value = 2.39036586939997
$ws1.Cells.Item(3, 29).Value = 0.732884
$ws1.Cells.Item(3, 30).Value = 0.266941
$ws1.Cells.Item(3, 31).Value = 3.74614615214598
$ws1.Cells.Item(3, 32).Value = 0.742102
$ws1.Cells.Item(3, 33).Value = 0.257898
$ws1.Cells.Item(3, 34).Value = 1.347523655777777
$ws1.Cells.Item(3, 35).Value = 0.494021
$ws1.Cells.Item(3, 36).Value = 0.505979
$ws1.Cells.Item(3, 37).Value = 2.024205448756227
$ws1.Cells.Item(3, 38).Value = 0.722187
$ws1.Cells.Item(3, 39).Value = 0.277813
$ws1.Cells.Item(3, 40).Value = 1.384682914535986
$ws1.Cells.Item(3, 41).Value = 0.468153
$ws1.Cells.Item(3, 42).Value = 0.531847
$ws1.Cells.Item(3, 43).Value = 2.136053811467618
$ws1.Cells.Item(3, 44).Value = 0.765397
$ws1.Cells.Item(3, 45).Value = 1.306511522778375
$ws1.Cells.Item(3, 46).Value = 0.736681
$ws1.Cells.Item(3, 47).Value = 1.357439651626688

# Row 4: Динамо М – Торпедо
$ws1.Cells.Item(4, 1).Value = 1369
$ws1.Cells.Item(4, 2).Value = 45992.8125
$ws1.Cells.Item(4, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(4, 3).Value = "Динамо М"
$ws1.Cells.Item(4, 4).Value = "Торпедо"
$ws1.Cells.Item(4, 5).Value = "Динамо М – Торпедо"
$ws1.Cells.Item(4, 6).Value = 897837
$ws1.Cells.Item(4, 7).Value = "https://text.khl.ru/text/897837.html"
$ws1.Cells.Item(4, 8).Value = 1.460695
$ws1.Cells.Item(4, 9).Value = 2.942447
$ws1.Cells.Item(4, 10).Value = 4.403142
$ws1.Cells.Item(4, 11).Value = 24.776157
$ws1.Cells.Item(4, 12).Value = 31.069775
$ws1.Cells.Item(4, 13).Value = 55.845933
$ws1.Cells.Item(4, 14).Value = 0.235408
$ws1.Cells.Item(4, 15).Value = 0.201686
$ws1.Cells.Item(4, 16).Value = 0.5628840000000001
$ws1.Cells.Item(4, 17).Value = 4.247943995106368
$ws1.Cells.Item(4, 18).Value = 4.958202354154477
$ws1.Cells.Item(4, 19).Value = 1.776564976087435
$ws1.Cells.Item(4, 20).Value = 23.5408
$ws1.Cells.Item(4, 21).Value = 20.1686
$ws1.Cells.Item(4, 22).Value = 56.2884
$ws1.Cells.Item(4, 23).Value = 0.703372
$ws1.Cells.Item(4, 24).Value = 0.296606
$ws1.Cells.Item(4, 25).Value = 3.371475964747848
$ws1.Cells.Item(4, 26).Value = 0.841905
$ws1.Cells.Item(4, 27).Value = 0.158073
$ws1.Cells.Item(4, 28).Value = 6.326191063622504
$ws1.Cells.Item(4, 29).Value = 0.925396
$ws1.Cells.Item(4, 30).Value = 0.074582
$ws1.Cells.Item(4, 31).Value = 13.40806092622885
$ws1.Cells.Item(4, 32).Value = 0.405801
$ws1.Cells.Item(4, 33).Value = 0.594199
$ws1.Cells.Item(4, 34).Value = 2.464262039768261
$ws1.Cells.Item(4, 35).Value = 0.16485
$ws1.Cells.Item(4, 36).Value = 0.8351499999999999
$ws1.Cells.Item(4, 37).Value = 6.066120715802245
$ws1.Cells.Item(4, 38).Value = 0.650971
$ws1.Cells.Item(4, 39).Value = 0.349029
$ws1.Cells.Item(4, 40).Value = 1.536166741682809
$ws1.Cells.Item(4, 41).Value = 0.383413
$ws1.Cells.Item(4, 42).Value = 0.616587
$ws1.Cells.Item(4, 43).Value = 2.608153609815004
$ws1.Cells.Item(4, 44).Value = 0.651632
$ws1.Cells.Item(4, 45).Value = 1.534608490681857
$ws1.Cells.Item(4, 46).Value = 0.899029
$ws1.Cells.Item(4, 47).Value = 1.112311171274786

# Row 5: ЦСКА – Салават Юлаев
$ws1.Cells.Item(5, 1).Value = 1369
$ws1.Cells.Item(5, 2).Value = 45992.8125
$ws1.Cells.Item(5, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws1.Cells.Item(5, 3).Value = "ЦСКА"
$ws1.Cells.Item(5, 4).Value = "Салават Юлаев"
$ws1.Cells.Item(5, 5).Value = "ЦСКА – Салават Юлаев"
$ws1.Cells.Item(5, 6).Value = 897834
$ws1.Cells.Item(5, 7).Value = "https://text.khl.ru/text/897834.html"
$ws1.Cells.Item(5, 8).Value = 3.561942
$ws1.Cells.Item(5, 9).Value = 1.027778
$ws1.Cells.Item(5, 10).Value = 4.58972
$ws1.Cells.Item(5, 11).Value = 29.27658
$ws1.Cells.Item(5, 12).Value = 21.790523
$ws1.Cells.Item(5, 13).Value = 51.067103
$ws1.Cells.Item(5, 14).Value = 0.736415
$ws1.Cells.Item(5, 15).Value = 0.146246
$ws1.Cells.Item(5, 16).Value = 0.117143
$ws1.Cells.Item(5, 17).Value = 1.357929971551367
$ws1.Cells.Item(5, 18).Value = 6.837793854190885
$ws1.Cells.Item(5, 19).Value = 8.536574955396397
$ws1.Cells.Item(5, 20).Value = 73.64150000000001
$ws1.Cells.Item(5, 21).Value = 14.6246
$ws1.Cells.Item(5, 22).Value = 11.7143
$ws1.Cells.Item(5, 23).Value = 0.637228
$ws1.Cells.Item(5, 24).Value = 0.362576
$ws1.Cells.Item(5, 25).Value = 2.758042451789418
$ws1.Cells.Item(5, 26).Value = 0.791807
$ws1.Cells.Item(5, 27).Value = 0.207997
$ws1.Cells.Item(5, 28).Value = 4.80776165040842
$ws1.Cells.Item(5, 29).Value = 0.893753
$ws1.Cells.Item(5, 30).Value = 0.106051
$ws1.Cells.Item(5, 31).Value = 9.42942546510641
$ws1.Cells.Item(5, 32).Value = 0.778914
$ws1.Cells.Item(5, 33).Value = 0.221086
$ws1.Cells.Item(5, 34).Value = 1.28383878066128
$ws1.Cells.Item(5, 35).Value = 0.5446839999999999
$ws1.Cells.Item(5, 36).Value = 0.455316
$ws1.Cells.Item(5, 37).Value = 1.83592688604769
$ws1.Cells.Item(5, 38).Value = 0.299953
$ws1.Cells.Item(5, 39).Value = 0.700047
$ws1.Cells.Item(5, 40).Value = 3.33385563738319
$ws1.Cells.Item(5, 41).Value = 0.099025
$ws1.Cells.Item(5, 42).Value = 0.900975
$ws1.Cells.Item(5, 43).Value = 10.09845998485231
$ws1.Cells.Item(5, 44).Value = 0.959053
$ws1.Cells.Item(5, 45).Value = 1.042695242077341
$ws1.Cells.Item(5, 46).Value = 0.462498
$ws1.Cells.Item(5, 47).Value = 2.162171512093025

# ---- Sheet2 (Cards_telegram): update rows 2-4, add row 5 ----
# Row 2: Сибирь – Барыс
$ws2.Cells.Item(2, 1).Value = 45992.64583333334
$ws2.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(2, 2).Value = "Сибирь – Барыс"
$cardText2 = @"
КХЛ • Регулярный чемпионат • 01.12.2025
Сибирь – Барыс
Ожидания модели (60’):
• Голы: λ_total ≈ 5.75 (2.90 : 2.85)
• Броски: SOG λ ≈ 53 (23 : 30)
Исход (60’), честные кф:
• П1: 42.3%  (Kмод 2.36)
• Х:  17.0%  (Kмод 5.87)
• П2: 40.6%  (Kмод 2.46)
Тоталы голов:
• ТМ 4.5: 31.9%  (Kмод 3.13)
• ТБ 4.5: 68.0%  (Kмод 1.47)
• ТМ 5.5: 48.6%  (Kмод 2.06)
• ТБ 5.5: 51.4%  (Kмод 1.95)
• ТМ 6.5: 64.6%  (Kмод 1.55)
• ТБ 6.5: 35.4%  (Kмод 2.83)
Индивидуальные тоталы:
• Сибирь ИТБ 1.5: 78.6% (Kмод 1.27)
• Сибирь ИТБ 2.5: 55.5% (Kмод 1.80)
• Барыс ИТБ 1.5: 77.7% (Kмод 1.29)
• Барыс ИТБ 2.5: 54.2% (Kмод 1.84)
Фора +1.5:
• Сибирь +1.5: 74.7% (Kмод 1.34)
• Барыс +1.5: 73.2% (Kмод 1.37)
"@
$ws2.Cells.Item(2, 3).Value = $cardText2

# Row 3: Северсталь – Трактор
$ws2.Cells.Item(3, 1).Value = 45992.79166666666
$ws2.Cells.Item(3, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(3, 2).Value = "Северсталь – Трактор"
$cardText3 = @"
КХЛ • Регулярный чемпионат • 01.12.2025
Северсталь – Трактор
Ожидания модели (60’):
• Голы: λ_total ≈ 5.20 (2.65 : 2.55)
• Броски: SOG λ ≈ 60 (24 : 36)
Исход (60’), честные кф:
• П1: 42.8%  (Kмод 2.34)
• Х:  18.0%  (Kмод 5.57)
• П2: 39.3%  (Kмод 2.55)
Тоталы голов:
• ТМ 4.5: 40.7%  (Kмод 2.46)
• ТБ 4.5: 59.3%  (Kмод 1.69)
• ТМ 5.5: 58.1%  (Kмод 1.72)
• ТБ 5.5: 41.8%  (Kмод 2.39)
• ТМ 6.5: 73.3%  (Kмод 1.36)
• ТБ 6.5: 26.7%  (Kмод 3.75)
Индивидуальные тоталы:
• Северсталь ИТБ 1.5: 74.2% (Kмод 1.35)
• Северсталь ИТБ 2.5: 49.4% (Kмод 2.02)
• Трактор ИТБ 1.5: 72.2% (Kмод 1.38)
• Трактор ИТБ 2.5: 46.8% (Kмод 2.14)
Фора +1.5:
• Северсталь +1.5: 76.5% (Kмод 1.31)
• Трактор +1.5: 73.7% (Kмод 1.36)
"@
$ws2.Cells.Item(3, 3).Value = $cardText3

# Row 4: Динамо М – Торпедо
$ws2.Cells.Item(4, 1).Value = 45992.8125
$ws2.Cells.Item(4, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(4, 2).Value = "Динамо М – Торпедо"
$cardText4 = @"
КХЛ • Регулярный чемпионат • 01.12.2025
Динамо М – Торпедо
Ожидания модели (60’):
• Голы: λ_total ≈ 3.62 (1.39 : 2.22)
• Броски: SOG λ ≈ 56 (25 : 31)
Исход (60’), честные кф:
• П1: 23.5%  (Kмод 4.25)
• Х:  20.2%  (Kмод 4.96)
• П2: 56.3%  (Kмод 1.78)
Тоталы голов:
• ТМ 4.5: 70.3%  (Kмод 1.42)
• ТБ 4.5: 29.7%  (Kмод 3.37)
• ТМ 5.5: 84.2%  (Kмод 1.19)
• ТБ 5.5: 15.8%  (Kмод 6.33)
• ТМ 6.5: 92.5%  (Kмод 1.08)
• ТБ 6.5: 7.5%  (Kмод 13.41)
Индивидуальные тоталы:
• Динамо М ИТБ 1.5: 40.6% (Kмод 2.46)
• Динамо М ИТБ 2.5: 16.5% (Kмод 6.07)
• Торпедо ИТБ 1.5: 65.1% (Kмод 1.54)
• Торпедо ИТБ 2.5: 38.3% (Kмод 2.61)
Фора +1.5:
• Динамо М +1.5: 65.2% (Kмод 1.53)
• Торпедо +1.5: 89.9% (Kмод 1.11)
"@
$ws2.Cells.Item(4, 3).Value = $cardText4

# Row 5: ЦСКА – Салават Юлаев
$ws2.Cells.Item(5, 1).Value = 45992.8125
$ws2.Cells.Item(5, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws2.Cells.Item(5, 2).Value = "ЦСКА – Салават Юлаев"
$cardText5 = @"
КХЛ • Регулярный чемпионат • 01.12.2025
ЦСКА – Салават Юлаев
Ожидания модели (60’):
• Голы: λ_total ≈ 3.96 (2.86 : 1.10)
• Броски: SOG λ ≈ 51 (29 : 22)
Исход (60’), честные кф:
• П1: 73.6%  (Kмод 1.36)
• Х:  14.6%  (Kмод 6.84)
• П2: 11.7%  (Kмод 8.54)
Тоталы голов:
• ТМ 4.5: 63.7%  (Kмод 1.57)
• ТБ 4.5: 36.3%  (Kмод 2.76)
• ТМ 5.5: 79.2%  (Kмод 1.26)
• ТБ 5.5: 20.8%  (Kмод 4.81)
• ТМ 6.5: 89.4%  (Kмод 1.12)
• ТБ 6.5: 10.6%  (Kмод 9.43)
Индивидуальные тоталы:
• ЦСКА ИТБ 1.5: 77.9% (Kмод 1.28)
• ЦСКА ИТБ 2.5: 54.5% (Kмод 1.84)
• Салават Юлаев ИТБ 1.5: 30.0% (Kмод 3.33)
• Салават Юлаев ИТБ 2.5: 9.9% (Kмод 10.10)
Фора +1.5:
• ЦСКА +1.5: 95.9% (Kмод 1.04)
• Салават Юлаев +1.5: 46.2% (Kмод 2.16)
"@
$ws2.Cells.Item(5, 3).Value = $cardText5

